# Updated symbol list on Fri Dec 23 05:24:04 UTC 2022 with GitHub Actions
# Applies the latest scraped price/volume refresh to the crypto tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value. Values must stay plain text (matching the
# original inline-string cells), so NumberFormat is forced to "@" before the
# assignment to stop Excel from auto-coercing numeric-looking strings into
# actual numbers.
$updates = [ordered]@{
    "D2"  = "246.42"
    "D3"  = "22.03"
    "D4"  = "5.430"
    "D5"  = "0.05785"
    "D6"  = "3.390"
    "D7"  = "6.324"
    "D8"  = "0.8181"
    "D9"  = "0.9527"
    "E9"  = "8FTXTokenFTTBestin24h"
    "D11" = "0.07506"
    "D12" = "0.03138"
    "D13" = "0.02999"
    "D14" = "4.176"
    "D15" = "0.09411"
    "D16" = "0.001602"
    "D17" = "0.04826"
    "D18" = "0.0005850"
    "D19" = "0.006187"
    "D20" = "0.004122"
    "D21" = "0.0009984"
    "D23" = "3.770"
    "D24" = "2.222"
    "D27" = "0.0003999"
    "D40" = "0.03893"
    "D41" = "0.006370"
    "D42" = "0.1075"
    "D43" = "0.003000"
    "E43" = "42CEJICEJI"
    "D44" = "0.006696"
    "D45" = "0.00005593"
    "D47" = "0.3800"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
